$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Washington State"
$ws.Range("A7").Value = "UCLA"
$ws.Range("B6").Value = "https://www.sports-reference.com/cbb/schools/washington-state/2023.html"
$ws.Range("B7").Value = "https://www.sports-reference.com/cbb/schools/ucla/2023.html"

$ws.Range("A8").Select()
